$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9: update title and link
$ws.Range("D9").Value = "MBA AI/BigData를 위한 수학, 통계학 교재"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/math-stat-for-mba-ai-bigdata/#utm_source=rss&utm_medium=rss&utm_campaign=math-stat-for-mba-ai-bigdata"

# Row 32: update title and link
$ws.Range("D32").Value = "Transfer learning /  Fine tuning"
$ws.Range("E32").Value = "https://dodonam.tistory.com/350"

# Row 36: update title and link
$ws.Range("D36").Value = "How to train your ViT? Data, Augmentation,  and Regularization in Vision Transformers"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/349"
